# Regenerated data from tools
# Updates the correlation-analysis numbers on the five per-tool sheets
# (all_tools, checker_framework, typestate_checker, infer, openjml) for
# rows 10-12 (num_snippets_warnings / num_warnings and the derived
# kendalls/spearmans stats), plus a couple of column-width tweaks on the
# "infer" sheet.

$wb = $excel.ActiveWorkbook

function Set-Row($ws, $row, $F, $G, $I, $J, $K, $L) {
    if ($F -ne $null) { $ws.Cells.Item($row, 6).Value = $F }
    if ($G -ne $null) { $ws.Cells.Item($row, 7).Value = $G }
    $ws.Cells.Item($row, 9).Value  = $I
    $ws.Cells.Item($row, 10).Value = $J
    $ws.Cells.Item($row, 11).Value = $K
    $ws.Cells.Item($row, 12).Value = $L
}

# ---- Sheet: all_tools ----
$ws = $wb.Worksheets.Item("all_tools")
Set-Row $ws 10 49 802 0.09989928307027053 0.3372677737988888 0.1559244610842745 0.2795654654406445
Set-Row $ws 11 49 802 0.08931863763595778 0.3695323356271473 0.1375884906942706 0.3406751919077503
Set-Row $ws 12 49 802 0.1193675355124785  0.2274845227161654 0.1890061392789798 0.1886589730863144

# ---- Sheet: checker_framework ----
$ws = $wb.Worksheets.Item("checker_framework")
Set-Row $ws 10 23 68 -0.1127953217534384 0.3246303816462205  -0.141785465198251  0.3260061663372782
Set-Row $ws 11 23 68 -0.1657843576414015 0.1297187088392248  -0.228149078380763   0.1110368177920856
Set-Row $ws 12 23 68 0.2135201189924139  0.04951827186158275 0.2585375793536707  0.06985387930445178

# ---- Sheet: typestate_checker ----
$ws = $wb.Worksheets.Item("typestate_checker")
Set-Row $ws 10 $null 520 0.08993343971074573 0.3939555332313259 0.1261078291007821  0.3828516225378101
Set-Row $ws 11 $null 520 0.07727463833790227 0.4435614692500816 0.09527592158164427 0.5104337983014946
Set-Row $ws 12 $null 520 0.06408768791953229 0.5225032783206667 0.08576989054836207 0.5536912247004897

# ---- Sheet: infer ----
$ws = $wb.Worksheets.Item("infer")
Set-Row $ws 10 2 2 -0.07463933708620761 0.547733910068501   -0.08587989564247843 0.5531808807861933
Set-Row $ws 11 2 2 -0.02642855544759036 0.8236209225496525  -0.03184146471615851 0.8262496889100787
Set-Row $ws 12 2 2 0.2391168558431198   0.04237662250330112 0.2899717804431688   0.04108542769241475

# Column-width tweaks on the "infer" sheet (J shrinks, K grows by the same amount)
$ws.Range("J1").ColumnWidth = 18.8
$ws.Range("K1").ColumnWidth = 20.8

# ---- Sheet: openjml ----
$ws = $wb.Worksheets.Item("openjml")
Set-Row $ws 10 $null 212 0.02781643858671354 0.7967876810602408 0.0356633742276846 0.8057766752165008
Set-Row $ws 11 $null 212 0.0464013402367485  0.6531397458483266 0.06459222671217832 0.6558502732479303
Set-Row $ws 12 $null 212 0.1651861681364787  0.1073776851625311 0.2294861708178067  0.1089023642733105
